# Remove the unnecessary "id_scenario" column from the worksheet and its
# backing Excel table (Table1). This is an interior column (column A), so it
# is removed by deleting the worksheet column (which shifts every other
# column left by one) and then rebuilding the table definition over the new
# A1:I9 range so the table's column list/ref stay in sync with the data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the table's cosmetic settings before we touch anything.
$lo = $ws.ListObjects.Item(1)
$tableName = $lo.Name
$tableStyle = "TableStyleMedium6"

# Delete column A (id_scenario) outright - the rest of the data shifts left.
$ws.Range("A1").EntireColumn.Delete()

# The old table definition still points at the stale A1:J9 / 10-column
# layout, so drop it and recreate it over the new A1:I9 range (with headers)
# so xl/tables/table1.xml reflects the correct columns.
$lo.Unlist()
$newTable = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $ws.Range("A1:I9"), $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$newTable.Name = $tableName
$newTable.TableStyle = $tableStyle

# Match the saved selection/cursor position recorded in the edited workbook.
$ws.Range("B11").Select() | Out-Null
